$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C (stateness), shifting stateness,
# successful_transition and approach one column to the right.
$ws.Range("C1").EntireColumn.Insert()

# New header for the inserted column; copy formatting from a neighboring
# header cell so it keeps the same bold/centered header style.
$ws.Range("D1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$ws.Range("C1").Value = "proxy_score"

# New proxy_score values (column C) and updated stateness values (column D)
$data = @(
    @(7.34216935367578,  91.77711692094725),
    @(3.475321885018063, 43.44152356272578),
    @(4.718378336635055, 58.97972920793819),
    @(2.189408382985685, 27.36760478732106),
    @(2.429787986552066, 30.37234983190082),
    @(2.273261133685839, 28.41576417107299),
    @(4.681216693802947, 58.51520867253684)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $data[$i][0]
    $ws.Cells.Item($row, 4).Value = $data[$i][1]
}
